$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell -> new text value. NumberFormat is forced to Text ("@") before
# assignment so Excel stores the literal string instead of auto-coercing it to
# a number/percentage, then the style is reset back to Normal so the cell keeps
# its original (default) formatting -- only the underlying text value changes.
$updates = @{
    "D2" = "290.78"
    "E2" = "-3.57%"
    "D3" = "30.78"
    "E3" = "-6.43%"
    "D4" = "4.958"
    "E4" = "-0.01%"
    "D5" = "0.07223"
    "E5" = "-7.13%"
    "D6" = "1.800"
    "E6" = "-8.26%"
    "D7" = "7.678"
    "E7" = "-2.24%"
    "D8" = "3.760"
    "E8" = "-0.99%"
    "D9" = "0.8957"
    "E9" = "-3.08%"
    "D10" = "0.1650"
    "E10" = "-6.59%"
    "D11" = "0.07681"
    "E11" = "-2.35%"
    "D12" = "0.08022"
    "E12" = "-7.23%"
    "D13" = "0.03043"
    "E13" = "-3.43%"
    "E14" = "0.10%"
    "D15" = "0.001507"
    "E15" = "-0.56%"
    "D16" = "0.005691"
    "E16" = "-4.06%"
    "D17" = "3.471"
    "E17" = "0.24%"
    "D18" = "2.084"
    "E18" = "-3.26%"
    "D19" = "0.3311"
    "E19" = "-0.80%"
    "D20" = "0.1310"
    "E20" = "-0.59%"
    "D21" = "4.038"
    "E21" = "-6.49%"
    "E22" = "9.68%"
    "D23" = "0.04513"
    "E23" = "-1.03%"
    "D24" = "0.001214"
    "E24" = "-0.91%"
    "E25" = "-9.54%"
    "E26" = "-0.07%"
    "E39" = "-6.33%"
    "D40" = "0.04398"
    "E40" = "-6.80%"
    "D41" = "0.007312"
    "E41" = "-6.80%"
    "D42" = "0.1309"
    "E42" = "-3.39%"
    "D43" = "0.007684"
    "D44" = "0.001901"
    "E44" = "-18.86%"
    "D45" = "0.009211"
    "E45" = "-12.76%"
    "D46" = "0.00005928"
    "E46" = "-5.37%"
    "E47" = "-0.07%"
    "D48" = "2.247"
    "E48" = "173.92%"
    "D49" = "0.003003"
    "E50" = "-0.07%"
    "E51" = "-0.07%"
}

foreach ($cell in $updates.Keys) {
    $range = $ws.Range($cell)
    $range.NumberFormat = "@"
    $range.Value = $updates[$cell]
    $range.Style = "Normal"
}
